# Added raw accuracy (and nan) column.
# Insert a new column before column C; Excel shifts BAARD2/BAARD3/FS/LID/RC
# blocks, the FPR/Success header row, and all data columns one slot to the
# right, inheriting the bold/bordered style from column B for the new column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C:C").Insert()

# The two empty placeholder cells that used to sit at the left edge of the
# BAARD2/BAARD3/FS/LID/RC header band (B1, B2) are not part of the new
# layout; drop them now that the header band has shifted into C:M.
$ws.Range("B1:B2").Clear()

# New header label for the inserted column.
$ws.Range("C3").Value = "Accuracy after attack"

# New "raw accuracy after attack" values for each (attack, epsilon) row.
$ws.Range("C4").Value = 82.79569892473118
$ws.Range("C5").Value = 5.376344086021505
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 93.54838709677421
$ws.Range("C9").Value = 91.39784946236558
$ws.Range("C10").Value = 81.72043010752688
$ws.Range("C11").Value = 49.46236559139785
$ws.Range("C12").Value = 12.90322580645161
$ws.Range("C13").Value = 84.94623655913979
$ws.Range("C14").Value = 59.13978494623656
$ws.Range("C15").Value = 4.301075268817205
$ws.Range("C16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("C18").Value = 49.46236559139785
$ws.Range("C19").Value = 7.526881720430108
$ws.Range("C20").Value = 7.526881720430108
$ws.Range("C21").Value = 7.526881720430108
$ws.Range("C22").Value = 1.075268817204301
$ws.Range("C23").Value = 7.526881720430108
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 82.79569892473118
$ws.Range("C26").Value = 6.451612903225806
$ws.Range("C27").Value = 0
$ws.Range("C28").Value = 0
